$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I (I0) and J (IF), matching the formatting of the
# existing header cells (bold, centered, bordered style reused from H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for the new I0 / IF columns, rows 2-36.
$data = @(
    @(9, 9),
    @(7, 8),
    @(6, 8),
    @(9, 9),
    @(2, 4),
    @(4, 5),
    @(5, 7),
    @(2, 4),
    @(8, 8),
    @(8, 9),
    @(10, 10),
    @(6, 7),
    @(8, 8),
    @(7, 7),
    @(7, 8),
    @(6, 6),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(6, 6),
    @(8, 9),
    @(7, 7),
    @(7, 7),
    @(9, 9),
    @(12, 12),
    @(9, 9),
    @(8, 9),
    @(8, 8),
    @(3, 4),
    @(8, 8),
    @(8, 8),
    @(3, 4),
    @(7, 7),
    @(5, 6),
    @(6, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
